$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NPC rows to append (id, name, real_name, type, game_map_id, x_position, y_position)
$newRows = @(
    @{ id = 40; name = "Ladyoftheflames"; real_name = "Lady of the flames"; type = 1; game_map_id = "Hell";                x = 1040; y = 1952 },
    @{ id = 41; name = "LadyoftheShade";  real_name = "Lady of the Shade";  type = 2; game_map_id = "Delusional Memories"; x = 544;  y = 224 },
    @{ id = 42; name = "EmeraldSoul";     real_name = "Emerald Soul";       type = 2; game_map_id = "Delusional Memories"; x = 1280; y = 2064 },
    @{ id = 43; name = "KeyMaker";        real_name = "Key Maker";         type = 2; game_map_id = "Labyrinth";           x = 1600; y = 208 }
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
if ($startRow -lt 2) { $startRow = 2 }

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.id
    $ws.Cells.Item($r, 2).Value = $row.name
    $ws.Cells.Item($r, 3).Value = $row.real_name
    $ws.Cells.Item($r, 4).Value = $row.type
    $ws.Cells.Item($r, 5).Value = $row.game_map_id
    $ws.Cells.Item($r, 9).Value = $row.x
    $ws.Cells.Item($r, 10).Value = $row.y
    $r++
}
